# Add a new "2022" data column (S) to the renewable-energy table, continuing
# the existing year-by-year series that currently ends at column R (2021).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for 2022, one row below the other existing year columns:
#   row 4 -> year header
#   row 5 -> share of renewable energy, in percent
#   row 6 -> hydropower electricity production, mln kWh
$ws.Range("S4").Value = 2022
$ws.Range("S5").Value = 30
$ws.Range("S6").Value = 11928.6

# Copy the formatting of the last existing data column (R) onto the new
# column (S) so the new cells keep the same borders/fonts/number formats
# as the rest of the table.
$ws.Range("R4:R6").Copy()
$ws.Range("S4:S6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# The author's saved selection after the edit was cell T3 (just past the
# new column), so mirror that here.
$ws.Range("T3").Select()
